$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# In the use-case table, the "アクター" (Actor) row's value cell currently
# reads "会員" (Member). It must become "受付・司書" (Receptionist/Librarian).
# The table has a single "アクター" label in column 1; its value lives in
# column 2 of the same row. We locate that row by scanning column 1, then
# overwrite the value cell's Range.Text directly (Find/Replace is avoided
# here because "会員" also appears verbatim elsewhere in the document, e.g.
# inside "会員台帳の会員情報の退会", and we only want this one cell touched).
$table = $d.Tables.Item(1)
$actorRowIndex = $null
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $row = $table.Rows.Item($r)
    if ($row.Cells.Count -ge 2) {
        $labelText = $row.Cells.Item(1).Range.Text
        if ($labelText.Substring(0, $labelText.Length - 2) -eq "アクター") {
            $actorRowIndex = $r
        }
    }
}

if ($actorRowIndex -ne $null) {
    $valueCell = $table.Rows.Item($actorRowIndex).Cells.Item(2)
    $valueCell.Range.Text = "受付・司書"
}

# --- Edit 2 -----------------------------------------------------------------
# A bullet paragraph currently holds three separate runs whose text
# concatenates to "システムは" + "退会" + "完了画面を表示する". They must be
# collapsed into a single run containing "システムは退会完了画面を表示する".
# Word's Find/Replace reads/writes the document's visible text stream
# (spanning run boundaries) and naturally merges the matched runs into a
# single run using the first matched run's formatting, which is exactly the
# desired result. The full phrase is unique in the document, so a plain
# whole-document Find/Replace (MatchWholeWord, no wrap) touches only this
# paragraph.
$d.Content.Find.Execute("システムは退会完了画面を表示する", $true, $false, $false, $false, $false, `
                         $true, 0, $false, "システムは退会完了画面を表示する", 1) | Out-Null
